# Update both data sheets (NBR and BAR) so that the sensitivity table only
# keeps the last 15 data rows of the original 19-row series (cutoffs shifted
# by 4), reflecting the re-run with the fixed workflow.

$wb = $excel.ActiveWorkbook

# New values: for each sheet, Column A keeps the 0-based index (0..14),
# column B and C are taken from the *old* rows 6..20 (i.e. shifted by 4).

$sheet1Data = @(
    @(0, 5, 816),
    @(1, 6, 814),
    @(2, 7, 813),
    @(3, 8, 799),
    @(4, 9, 810),
    @(5, 10, 805),
    @(6, 11, 798),
    @(7, 12, 787),
    @(8, 13, 788),
    @(9, 14, 790),
    @(10, 15, 783),
    @(11, 16, 776),
    @(12, 17, 773),
    @(13, 18, 772),
    @(14, 19, 770)
)

$sheet2Data = @(
    @(0, 5, 1209),
    @(1, 6, 1212),
    @(2, 7, 1212),
    @(3, 8, 1199),
    @(4, 9, 1176),
    @(5, 10, 1181),
    @(6, 11, 1184),
    @(7, 12, 1195),
    @(8, 13, 1198),
    @(9, 14, 1191),
    @(10, 15, 1197),
    @(11, 16, 1201),
    @(12, 17, 1200),
    @(13, 18, 1197),
    @(14, 19, 1201)
)

$sheetData = @{
    "NBR" = $sheet1Data
    "BAR" = $sheet2Data
}

foreach ($sheetName in $sheetData.Keys) {
    $ws = $wb.Worksheets.Item($sheetName)
    $data = $sheetData[$sheetName]

    # Remove the now-unused trailing rows (old rows 17-20) so the sheet
    # dimension shrinks from A1:C20 down to A1:C16.
    $ws.Range("A17:C20").EntireRow.Delete() | Out-Null

    for ($i = 0; $i -lt $data.Length; $i++) {
        $row = 2 + $i
        $vals = $data[$i]
        $ws.Cells.Item($row, 1).Value = $vals[0]
        $ws.Cells.Item($row, 2).Value = $vals[1]
        $ws.Cells.Item($row, 3).Value = $vals[2]
    }
}
